# Update "想去人数" (column F) counts on sheets "展览" and "全部类型".
# The mapping below is keyed by exact row number on each sheet (safer
# than matching by old value, since some values repeat across rows).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 191
    4  = 183
    5  = 5146
    9  = 572
    10 = 526
    14 = 4100
    16 = 161
    17 = 143
    19 = 3094
    21 = 1045
    22 = 96
    24 = 185
    25 = 87
    26 = 21
    30 = 11
    33 = 7
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 191
    4  = 183
    6  = 5146
    10 = 572
    11 = 526
    15 = 4100
    17 = 161
    18 = 143
    20 = 3094
    22 = 1045
    23 = 96
    25 = 185
    26 = 87
    27 = 21
    31 = 11
    34 = 7
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
